# Applies the betting-odds value updates for 2025-11-24 (Betfair Back/Lay workbook).
# Generated from the authoritative cell-level diff; each assignment sets a single
# numeric cell to its new value while leaving every other cell untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 5.5
$ws.Range("I2").Value = 7.6
$ws.Range("J2").Value = 4
$ws.Range("K2").Value = 4.7
$ws.Range("V2").Value = 1.16

# Row 3
$ws.Range("F3").Value = 2.18
$ws.Range("J3").Value = 3.25
$ws.Range("K3").Value = 3.75
$ws.Range("L3").Value = 1.4
$ws.Range("Z3").Value = 27
$ws.Range("AJ3").Value = 34
$ws.Range("AL3").Value = 44
$ws.Range("AN3").Value = 22

# Row 4
$ws.Range("F4").Value = 2.58
$ws.Range("H4").Value = 3.05
$ws.Range("I4").Value = 3.4
$ws.Range("J4").Value = 2.94
$ws.Range("K4").Value = 3.35
$ws.Range("L4").Value = 1.53
$ws.Range("M4").Value = 1.11
$ws.Range("N4").Value = 2.72
$ws.Range("O4").Value = 1.49
$ws.Range("P4").Value = 1.58
$ws.Range("Q4").Value = 2.44
$ws.Range("R4").Value = 1.22
$ws.Range("S4").Value = 4.8
$ws.Range("T4").Value = 1.99
$ws.Range("U4").Value = 1.83
$ws.Range("V4").Value = 1.42
$ws.Range("X4").Value = 11.5
$ws.Range("Y4").Value = 11.5
$ws.Range("AB4").Value = 10.5
$ws.Range("AC4").Value = 8.6
$ws.Range("AH4").Value = 26
$ws.Range("AK4").Value = 44
$ws.Range("AL4").Value = 65
$ws.Range("AM4").Value = 190
$ws.Range("AO4").Value = 60

# Row 5
$ws.Range("F5").Value = 3.55
$ws.Range("G5").Value = 3.95
$ws.Range("H5").Value = 2.2
$ws.Range("I5").Value = 2.42
$ws.Range("J5").Value = 3.15
$ws.Range("L5").Value = 1.47
$ws.Range("M5").Value = 1.09
$ws.Range("Q5").Value = 2.16
$ws.Range("T5").Value = 1.87
$ws.Range("U5").Value = 1.84
$ws.Range("W5").Value = 1.34
$ws.Range("Y5").Value = 980
$ws.Range("AJ5").Value = 1000
$ws.Range("AK5").Value = 60
$ws.Range("AL5").Value = 65
$ws.Range("AN5").Value = 60

# Row 6
$ws.Range("F6").Value = 3.45
$ws.Range("J6").Value = 3.25
$ws.Range("L6").Value = 1.49
$ws.Range("M6").Value = 1.1
$ws.Range("Q6").Value = 2.26
$ws.Range("T6").Value = 1.92
$ws.Range("X6").Value = 13

# Row 7
$ws.Range("F7").Value = 1.67
$ws.Range("G7").Value = 1.79
$ws.Range("H7").Value = 4.8
$ws.Range("K7").Value = 4.6
$ws.Range("M7").Value = 1.05
$ws.Range("N7").Value = 4.2
$ws.Range("P7").Value = 2.18
$ws.Range("R7").Value = 1.44
$ws.Range("T7").Value = 1.73
$ws.Range("U7").Value = 2.08
$ws.Range("W7").Value = 2.26
$ws.Range("Z7").Value = 46
$ws.Range("AA7").Value = 160
$ws.Range("AC7").Value = 10.5
$ws.Range("AE7").Value = 80
$ws.Range("AG7").Value = 11

# Row 8
$ws.Range("F8").Value = 1.87
$ws.Range("M8").Value = 1.07
$ws.Range("O8").Value = 1.32

# Row 9
$ws.Range("H9").Value = 2.38
$ws.Range("I9").Value = 2.64
$ws.Range("Q9").Value = 1.76
$ws.Range("R9").Value = 1.43
$ws.Range("U9").Value = 2.3
$ws.Range("AJ9").Value = 980

# Row 10
$ws.Range("F10").Value = 1.8
$ws.Range("G10").Value = 1.84
$ws.Range("H10").Value = 5.2
$ws.Range("O10").Value = 1.36
$ws.Range("Q10").Value = 2.06
$ws.Range("R10").Value = 1.32
$ws.Range("T10").Value = 1.94
$ws.Range("W10").Value = 2.18
$ws.Range("AH10").Value = 980

# Row 11
$ws.Range("F11").Value = 2.58
$ws.Range("G11").Value = 2.7
$ws.Range("I11").Value = 3.05
$ws.Range("T11").Value = 1.68
$ws.Range("V11").Value = 1.48
$ws.Range("W11").Value = 1.58

# Row 12
$ws.Range("F12").Value = 3.6
$ws.Range("G12").Value = 3.65
$ws.Range("H12").Value = 2.38
$ws.Range("I12").Value = 2.42
$ws.Range("L12").Value = 1.55
$ws.Range("M12").Value = 1.11
$ws.Range("N12").Value = 2.96
$ws.Range("O12").Value = 1.49
$ws.Range("P12").Value = 1.64
$ws.Range("R12").Value = 1.23
$ws.Range("V12").Value = 1.7
$ws.Range("W12").Value = 1.37
$ws.Range("X12").Value = 9.199999999999999
$ws.Range("Y12").Value = 8.199999999999999
$ws.Range("Z12").Value = 13
$ws.Range("AA12").Value = 32
$ws.Range("AB12").Value = 11
$ws.Range("AE12").Value = 30
$ws.Range("AF12").Value = 22
$ws.Range("AG12").Value = 15.5
$ws.Range("AH12").Value = 22
$ws.Range("AJ12").Value = 70
$ws.Range("AO12").Value = 29

# Row 13
$ws.Range("P13").Value = 2.24
$ws.Range("R13").Value = 1.5
$ws.Range("T13").Value = 1.62
$ws.Range("U13").Value = 2.4
$ws.Range("AK13").Value = 28

# Row 14
$ws.Range("F14").Value = 2.34
$ws.Range("V14").Value = 1.41
$ws.Range("AN14").Value = 25

# Row 15
$ws.Range("L15").Value = 1.43
$ws.Range("N15").Value = 3.35
$ws.Range("P15").Value = 1.82
$ws.Range("R15").Value = 1.31

# Row 16
$ws.Range("N16").Value = 4.9

# Row 17
$ws.Range("F17").Value = 1.97
$ws.Range("I17").Value = 3.95
$ws.Range("J17").Value = 3.9
$ws.Range("K17").Value = 4.6
$ws.Range("P17").Value = 2.36
$ws.Range("R17").Value = 1.54
$ws.Range("U17").Value = 2.42
$ws.Range("V17").Value = 1.34

# Row 19
$ws.Range("AG19").Value = 11.5

# Row 20
$ws.Range("F20").Value = 2.86
$ws.Range("I20").Value = 3.1
$ws.Range("L20").Value = 1.47
$ws.Range("M20").Value = 1.13
$ws.Range("N20").Value = 2.68
$ws.Range("Q20").Value = 2.6
$ws.Range("R20").Value = 1.21
$ws.Range("V20").Value = 1.47
$ws.Range("Y20").Value = 9
$ws.Range("AD20").Value = 14.5
$ws.Range("AH20").Value = 25
$ws.Range("AI20").Value = 70

# Row 21
$ws.Range("F21").Value = 2.54
$ws.Range("K21").Value = 3.35
$ws.Range("L21").Value = 1.51
$ws.Range("AI21").Value = 70
$ws.Range("AO21").Value = 65

# Row 22
$ws.Range("K22").Value = 3.45
$ws.Range("P22").Value = 1.72
$ws.Range("Q22").Value = 2.22
$ws.Range("R22").Value = 1.27
$ws.Range("W22").Value = 1.3

# Row 23
$ws.Range("H23").Value = 4.3
$ws.Range("Q23").Value = 2.38

# Row 24
$ws.Range("F24").Value = 2.26
$ws.Range("G24").Value = 2.48
$ws.Range("J24").Value = 2.78
$ws.Range("K24").Value = 2.98
$ws.Range("W24").Value = 1.69

# Row 25
$ws.Range("S25").Value = 2.98
$ws.Range("T25").Value = 1.77
$ws.Range("U25").Value = 2.24
$ws.Range("AD25").Value = 20
$ws.Range("AH25").Value = 19
$ws.Range("AJ25").Value = 18.5
$ws.Range("AM25").Value = 85
$ws.Range("AN25").Value = 9.4

# Row 26
$ws.Range("P26").Value = 1.77

# Row 27
$ws.Range("F27").Value = 1.65
$ws.Range("J27").Value = 3.75

# Row 29
$ws.Range("Q29").Value = 2.3
$ws.Range("AA29").Value = 50

# Row 30
$ws.Range("F30").Value = 2.16
$ws.Range("G30").Value = 2.32
$ws.Range("H30").Value = 3.55
$ws.Range("K30").Value = 3.6
$ws.Range("L30").Value = 1.42
$ws.Range("Q30").Value = 1.98
$ws.Range("T30").Value = 1.79
$ws.Range("X30").Value = 16.5

# Row 31
$ws.Range("G31").Value = 2.28
$ws.Range("I31").Value = 3.9
$ws.Range("P31").Value = 1.75
$ws.Range("W31").Value = 1.78
$ws.Range("AG31").Value = 14

# Row 32
$ws.Range("U32").Value = 1.64
